# Auto-generated: update D (Price) and E (Volume 1h) columns for the cryptos sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") '75.874.59'
Set-TextValue $ws.Range("E2") '  +1.59%  '
Set-TextValue $ws.Range("D3") '2.903.51'
Set-TextValue $ws.Range("E3") '  +2.44%  '
Set-TextValue $ws.Range("E4") '  +0.05%  '
Set-TextValue $ws.Range("D5") '197.77'
Set-TextValue $ws.Range("E5") '  +4.89%  '
Set-TextValue $ws.Range("D6") '595.50'
Set-TextValue $ws.Range("E6") '  -0.82%  '
Set-TextValue $ws.Range("E7") '  -0.01%  '
Set-TextValue $ws.Range("E8") '  -1.36%  '
Set-TextValue $ws.Range("D9") '0.194'
Set-TextValue $ws.Range("E9") '  +0.89%  '
Set-TextValue $ws.Range("D10") '2.902.31'
Set-TextValue $ws.Range("E10") '  +2.46%  '
Set-TextValue $ws.Range("D11") '0.421'
Set-TextValue $ws.Range("E11") '  +13.39%  '
Set-TextValue $ws.Range("E12") '  -0.98%  '
Set-TextValue $ws.Range("D13") '4.87'
Set-TextValue $ws.Range("E13") '  -0.59%  '
Set-TextValue $ws.Range("D14") '3.440.95'
Set-TextValue $ws.Range("D15") '75.787.30'
Set-TextValue $ws.Range("E15") '  +0.99%  '
Set-TextValue $ws.Range("D16") '0.0000188'
Set-TextValue $ws.Range("E16") '  +0.22%  '
Set-TextValue $ws.Range("D17") '27.25'
Set-TextValue $ws.Range("E17") '  +0.05%  '
Set-TextValue $ws.Range("D18") '2.905.81'
Set-TextValue $ws.Range("E18") '  +2.35%  '
Set-TextValue $ws.Range("D19") '8.80'
Set-TextValue $ws.Range("E19") '  -3.73%  '
Set-TextValue $ws.Range("D20") '12.70'
Set-TextValue $ws.Range("E20") '  +2.23%  '
Set-TextValue $ws.Range("D21") '376.11'
Set-TextValue $ws.Range("E21") '  +0.04%  '
Set-TextValue $ws.Range("D22") '2.29'
Set-TextValue $ws.Range("E22") '  +0.97%  '
Set-TextValue $ws.Range("E23") '  +0.97%  '
Set-TextValue $ws.Range("D24") '71.19'
Set-TextValue $ws.Range("E24") '  +0.67%  '
Set-TextValue $ws.Range("D25") '0.999'
Set-TextValue $ws.Range("E25") '  -0.06%  '
Set-TextValue $ws.Range("D26") '3.055.55'
Set-TextValue $ws.Range("E26") '  +2.24%  '
Set-TextValue $ws.Range("E27") '  -0.94%  '
Set-TextValue $ws.Range("D28") '9.59'
Set-TextValue $ws.Range("E28") '  -0.19%  '
Set-TextValue $ws.Range("D29") '0.0000108'
Set-TextValue $ws.Range("E29") '  +4.43%  '
Set-TextValue $ws.Range("E30") '  -0.11%  '
Set-TextValue $ws.Range("E31") '  -0.84%  '
Set-TextValue $ws.Range("D32") '500.09'
Set-TextValue $ws.Range("E32") '  -4.94%  '
Set-TextValue $ws.Range("D33") '7.68'
Set-TextValue $ws.Range("E33") '  -3.15%  '
Set-TextValue $ws.Range("D34") '1.79'
Set-TextValue $ws.Range("E34") '  -0.87%  '
Set-TextValue $ws.Range("E35") '  +0.10%  '
Set-TextValue $ws.Range("D36") '164.06'
Set-TextValue $ws.Range("E36") '  +1.19%  '
Set-TextValue $ws.Range("D37") '20.08'
Set-TextValue $ws.Range("E37") '  +0.02%  '
Set-TextValue $ws.Range("D39") '0.112'
Set-TextValue $ws.Range("E39") '  -6.25%  '
Set-TextValue $ws.Range("E40") '  -0.06%  '
Set-TextValue $ws.Range("D41") '179.56'
Set-TextValue $ws.Range("E41") '  -1.48%  '
Set-TextValue $ws.Range("D42") '0.342'
Set-TextValue $ws.Range("E42") '  +0.53%  '
Set-TextValue $ws.Range("E44") '  -2.63%  '
Set-TextValue $ws.Range("D45") '0.0909'
Set-TextValue $ws.Range("E45") '  +7.05%  '
Set-TextValue $ws.Range("D46") '40.12'
Set-TextValue $ws.Range("E46") '  +1.39%  '
Set-TextValue $ws.Range("E47") '  -3.72%  '
Set-TextValue $ws.Range("D48") '2.30'
Set-TextValue $ws.Range("E48") '  -2.95%  '
Set-TextValue $ws.Range("D49") '0.574'
Set-TextValue $ws.Range("E49") '  +0.68%  '
Set-TextValue $ws.Range("E50") '  +7.45%  '
Set-TextValue $ws.Range("E51") '  -1.20%  '
